$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Set the new actual-time value for Pandas (row 8, column D)
$ws.Range("D8").Value = 5.5

# Update the selection / view to reflect where the user ended up editing
$excel.ActiveWindow.ScrollRow = 2
$ws.Range("D14").Select()
